$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B30 was previously entered as text "4"; normalize it to a real number.
$ws.Range("B30").Value = 4

# Append a new annotation row (row 31) for Ruilin.
$ws.Range("A31").Value = "Ruilin"

# B31 must stay a text value "3" (not a number). Use a leading apostrophe to
# force text entry, then strip the quote-prefix formatting it introduces so
# the cell keeps the sheet's default (unstyled) look.
$ws.Range("B31").Value = "'3"
$ws.Range("B31").ClearFormats()

$ws.Range("C31").Value = "无"
$ws.Range("D31").Value = "DFT"
$ws.Range("E31").Value = "WRI"
$ws.Range("F31").Value = "74483628-1e12-4bb7-acfc-2ccaf38e6d81"
$ws.Range("G31").Value = "HyIFzx-0b_annotated.xlsx"
$ws.Range("H31").Value = "Acronyms are not properly defined."
